$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 data
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44911
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100103
$ws.Range("H49").Value = "Frutos de hueso (carozo)"
$ws.Range("I49").Value = 100103003
$ws.Range("J49").Value = "Damasco"
$ws.Range("K49").Value = "Dina"
$ws.Range("L49").Value = "Especial"
$ws.Range("M49").Value = 250
$ws.Range("N49").Value = 20000
$ws.Range("O49").Value = 20000
$ws.Range("P49").Value = 20000
$ws.Range("Q49").Value = "`$/caja 16 kilos"
$ws.Range("R49").Value = "Región de O'Higgins"
$ws.Range("S49").Value = 1250
$ws.Range("T49").Value = 16

# Row 50 data
$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 44911
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103003
$ws.Range("J50").Value = "Damasco"
$ws.Range("K50").Value = "Dina"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 200
$ws.Range("N50").Value = 18000
$ws.Range("O50").Value = 18000
$ws.Range("P50").Value = 18000
$ws.Range("Q50").Value = "`$/caja 16 kilos"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 1125
$ws.Range("T50").Value = 16

# Copy the date style from D48 to D49:D50
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D49:D50").PasteSpecial(-4122) | Out-Null
